# Update '想去人数' (interest count) values in column F across all sheets
# per gh-pages data refresh generated at 456a3b4
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F7").Value = 5927
$ws.Range("F8").Value = 10084
$ws.Range("F10").Value = 72
$ws.Range("F12").Value = 4000
$ws.Range("F13").Value = 218
$ws.Range("F14").Value = 145
$ws.Range("F16").Value = 122
$ws.Range("F18").Value = 685
$ws.Range("F19").Value = 3982
$ws.Range("F20").Value = 148
$ws.Range("F22").Value = 5567
$ws.Range("F23").Value = 448
$ws.Range("F24").Value = 2195
$ws.Range("F25").Value = 146
$ws.Range("F26").Value = 392
$ws.Range("F27").Value = 8325
$ws.Range("F30").Value = 2233
$ws.Range("F31").Value = 2266
$ws.Range("F32").Value = 1351
$ws.Range("F33").Value = 191
$ws.Range("F34").Value = 1470
$ws.Range("F35").Value = 25
$ws.Range("F36").Value = 293
$ws.Range("F41").Value = 1201
$ws.Range("F42").Value = 1193
$ws.Range("F46").Value = 1382
$ws.Range("F47").Value = 243

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F13").Value = 133

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 624
$ws.Range("F3").Value = 818
$ws.Range("F4").Value = 79

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F4").Value = 624
$ws.Range("F5").Value = 818
$ws.Range("F6").Value = 79
$ws.Range("F8").Value = 5927
$ws.Range("F9").Value = 10084
$ws.Range("F11").Value = 4000
$ws.Range("F12").Value = 218
$ws.Range("F14").Value = 122
$ws.Range("F18").Value = 685
$ws.Range("F19").Value = 3982
$ws.Range("F21").Value = 148
$ws.Range("F22").Value = 5567
$ws.Range("F23").Value = 448
$ws.Range("F24").Value = 2195
$ws.Range("F25").Value = 146
$ws.Range("F26").Value = 392
$ws.Range("F27").Value = 8325
$ws.Range("F30").Value = 2233
$ws.Range("F31").Value = 2266
$ws.Range("F32").Value = 1351
$ws.Range("F33").Value = 191
$ws.Range("F34").Value = 1472
$ws.Range("F35").Value = 25
$ws.Range("F36").Value = 293
$ws.Range("F40").Value = 1201
$ws.Range("F41").Value = 1193
$ws.Range("F45").Value = 1382
$ws.Range("F48").Value = 243

